$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Jason Lopez"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "109.64"

# Row 3
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "115.23"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "112346270283492312"
$ws.Range("E3").Value = "jason.lopez+2@tribal.credit"

# Row 4
$ws.Range("A4").Value = "Melissa Espinoza"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "112346270283492312"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "768"

# Row 5
$ws.Range("B5").Value = "Jason Lopez"
$ws.Range("F5").Value = "Testing 2"
